$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '34.437.82'
$ws.Range('E2').Value = '  +0.21%  '
$ws.Range('D3').Value = '1.804.02'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '224.66'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E6').Value = '  +2.82%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '38.25'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +5.75%  '
$ws.Range('E9').Value = '  -4.56%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0671'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.31%  '
$ws.Range('E11').Value = '  +0.97%  '
$ws.Range('D12').Value = '2.064.64'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.05'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -6.41%  '
$ws.Range('D14').Value = '1.807.30'
$ws.Range('E14').Value = '  +0.52%  '
$ws.Range('D15').Value = '34.423.12'
$ws.Range('E15').Value = '  +0.23%  '
$ws.Range('E16').Value = '  -2.79%  '
$ws.Range('E17').Value = '  -3.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '67.81'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.02%  '
$ws.Range('E19').Value = '  -1.68%  '
$ws.Range('D20').Value = '0.0₃0767'
$ws.Range('E20').Value = '  -3.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.05'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.52%  '
$ws.Range('E23').Value = '  -1.94%  '
$ws.Range('E24').Value = '  +2.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '170.81'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.37%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.68'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.39'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.05%  '
$ws.Range('E28').Value = '  +0.90%  '
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('E30').Value = '  -1.35%  '
$ws.Range('E31').Value = '  -2.32%  '
$ws.Range('E32').Value = '  -3.93%  '
$ws.Range('E34').Value = '  -0.38%  '
$ws.Range('D35').Value = '1.322.07'
$ws.Range('E35').Value = '  -5.46%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.639'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.94%  '
$ws.Range('E37').Value = '  -1.64%  '
$ws.Range('E38').Value = '  -1.35%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.30'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.30%  '
$ws.Range('B40').Value = 'HuobiToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.44'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.84%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '82.34'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.30%  '
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.81'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.73%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.21'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.95%  '
$ws.Range('E44').Value = '  -2.41%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.67'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.86%  '
$ws.Range('E46').Value = '  +0.40%  '
$ws.Range('D47').Value = '1.965.47'
$ws.Range('E47').Value = '  +0.12%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.74'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.66%  '
$ws.Range('E49').Value = '  -0.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '102.05'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.34%  '
$ws.Range('D51').Value = '0.0₆0119'
$ws.Range('E51').Value = '  -7.23%  '
